$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header in F1, matching the style used by other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:52:43.383504",
    "2021-10-05 10:52:43.383516",
    "2021-10-05 10:52:43.383519",
    "2021-10-05 10:52:43.383522",
    "2021-10-05 10:52:43.383525",
    "2021-10-05 10:52:43.383528",
    "2021-10-05 10:52:43.383530",
    "2021-10-05 10:52:43.383533",
    "2021-10-05 10:52:43.383535",
    "2021-10-05 10:52:43.383538",
    "2021-10-05 10:52:43.383541",
    "2021-10-05 10:52:43.383544",
    "2021-10-05 10:52:43.383546",
    "2021-10-05 10:52:43.383549",
    "2021-10-05 10:52:43.383551",
    "2021-10-05 10:52:43.383554",
    "2021-10-05 10:52:43.383556",
    "2021-10-05 10:52:43.383559",
    "2021-10-05 10:52:43.383562",
    "2021-10-05 10:52:43.383564",
    "2021-10-05 10:52:43.383567",
    "2021-10-05 10:52:43.383569",
    "2021-10-05 10:52:43.383572",
    "2021-10-05 10:52:43.383574",
    "2021-10-05 10:52:43.383577",
    "2021-10-05 10:52:43.383580",
    "2021-10-05 10:52:43.383582",
    "2021-10-05 10:52:43.383585",
    "2021-10-05 10:52:43.383587",
    "2021-10-05 10:52:43.383590",
    "2021-10-05 10:52:43.383592",
    "2021-10-05 10:52:43.383595",
    "2021-10-05 10:52:43.383598",
    "2021-10-05 10:52:43.383601"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
